$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 344.94736
$ws.Range("I28").Value = 280.53333
$ws.Range("J28").Value = 586.5
$ws.Range("K28").Value = 280.53333
$ws.Range("L28").Value = 586.5
$ws.Range("M28").Value = 204.46667
$ws.Range("N28").Value = -1556.5
$ws.Range("H88").Value = 604.38464
$ws.Range("I88").Value = 436.25
$ws.Range("J88").Value = 679.1111
$ws.Range("K88").Value = 436.25
$ws.Range("L88").Value = 679.1111
$ws.Range("M88").Value = -30.25
$ws.Range("N88").Value = -1491.1111
$ws.Range("H91").Value = 604.38464
$ws.Range("I91").Value = 436.25
$ws.Range("J91").Value = 679.1111
$ws.Range("K91").Value = 436.25
$ws.Range("L91").Value = 679.1111
$ws.Range("M91").Value = 967.75
$ws.Range("N91").Value = -3487.1111
$ws.Range("H113").Value = 58827052
$ws.Range("I113").Value = 142858700
$ws.Range("J113").Value = 4897.9
$ws.Range("K113").Value = 142858700
$ws.Range("L113").Value = 4897.9
$ws.Range("M113").Value = -142855446
$ws.Range("N113").Value = -11405.9
$ws.Range("H116").Value = 16671140
$ws.Range("I116").Value = 50000960
$ws.Range("J116").Value = 6229.4
$ws.Range("K116").Value = 50000960
$ws.Range("L116").Value = 6229.4
$ws.Range("M116").Value = -49997518
$ws.Range("N116").Value = -13113.4
$ws.Range("H127").Value = 979.4666999999999
$ws.Range("I127").Value = 470.85715
$ws.Range("J127").Value = 1424.5
$ws.Range("K127").Value = 1412.57145
$ws.Range("L127").Value = 4273.5
$ws.Range("M127").Value = 3547.42855
$ws.Range("N127").Value = -14193.5
$ws.Range("H129").Value = 159764.7
$ws.Range("J129").Value = 179698.17
$ws.Range("L129").Value = 539094.51
$ws.Range("N129").Value = -549094.51
$ws.Range("H131").Value = 2206.2727
$ws.Range("I131").Value = 1203.9
$ws.Range("J131").Value = 3041.5833
$ws.Range("K131").Value = 3611.7
$ws.Range("L131").Value = 9124.749899999999
$ws.Range("M131").Value = 1428.3
$ws.Range("N131").Value = -19204.7499
$ws.Range("H132").Value = 2411.1277
$ws.Range("I132").Value = 2420.575
$ws.Range("K132").Value = 7261.724999999999
$ws.Range("M132").Value = -4731.724999999999
$ws.Range("H141").Value = 3719.889
$ws.Range("I141").Value = 3679
$ws.Range("J141").Value = 3801.6667
$ws.Range("K141").Value = 11037
$ws.Range("L141").Value = 11405.0001
$ws.Range("M141").Value = -5857
$ws.Range("N141").Value = -21765.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H32").Value = 5191.924
$ws.Range("I32").Value = 3854.6904
$ws.Range("J32").Value = 19232.875
$ws.Range("K32").Value = 3854.6904
$ws.Range("L32").Value = 19232.875
$ws.Range("M32").Value = -3567.6904
$ws.Range("N32").Value = -19806.875
$ws.Range("H122").Value = 1791.7916
$ws.Range("I122").Value = 1723.9524
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 5171.857199999999
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -2721.857199999999
$ws.Range("N122").Value = -11700.0001
$ws.Range("H132").Value = 9488.462
$ws.Range("I132").Value = 1747.2979
$ws.Range("K132").Value = 5241.893700000001
$ws.Range("M132").Value = -2711.893700000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1909.5312
$ws.Range("I86").Value = 1771.6666
$ws.Range("K86").Value = 1771.6666
$ws.Range("M86").Value = -648.6666
$ws.Range("H89").Value = 1909.5312
$ws.Range("I89").Value = 1771.6666
$ws.Range("K89").Value = 8858.333000000001
$ws.Range("M89").Value = -3242.333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("H62").Value = 5833.3335
$ws.Range("I62").Value = 5833.3335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5833.3335
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5209.3335
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5833.3335
$ws.Range("I65").Value = 5833.3335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 29166.6675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -26046.6675
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 54845
$ws.Range("J68").Value = 54845
$ws.Range("L68").Value = 54845
$ws.Range("N68").Value = -56343
$ws.Range("H71").Value = 54845
$ws.Range("J71").Value = 54845
$ws.Range("L71").Value = 164535
$ws.Range("N71").Value = -172023

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4005.6924
$ws.Range("I3").Value = 1775.5555
$ws.Range("J3").Value = 9023.5
$ws.Range("K3").Value = 5326.666499999999
$ws.Range("L3").Value = 27070.5
$ws.Range("M3").Value = -5214.666499999999
$ws.Range("N3").Value = -27294.5
$ws.Range("H115").Value = 5002.4
$ws.Range("I115").Value = 30
$ws.Range("J115").Value = 5554.8887
$ws.Range("K115").Value = 90
$ws.Range("L115").Value = 16664.6661
$ws.Range("M115").Value = 1085
$ws.Range("N115").Value = -19014.6661
$ws.Range("H131").Value = 704.4400000000001
$ws.Range("J131").Value = 721.01044
$ws.Range("L131").Value = 2163.03132
$ws.Range("N131").Value = -12243.03132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3452.5
$ws.Range("I80").Value = 3012.3076
$ws.Range("K80").Value = 3012.3076
$ws.Range("M80").Value = -2014.3076
$ws.Range("H83").Value = 3452.5
$ws.Range("I83").Value = 3012.3076
$ws.Range("K83").Value = 15061.538
$ws.Range("M83").Value = -10069.538
$ws.Range("H113").Value = 4712.385
$ws.Range("I113").Value = 5957.9443
$ws.Range("J113").Value = 1909.875
$ws.Range("K113").Value = 5957.9443
$ws.Range("L113").Value = 1909.875
$ws.Range("M113").Value = -3787.9443
$ws.Range("N113").Value = -6249.875
$ws.Range("H126").Value = 3956.1555
$ws.Range("J126").Value = 3511.7144
$ws.Range("L126").Value = 10535.1432
$ws.Range("N126").Value = -15475.1432
$ws.Range("H132").Value = 30522.8
$ws.Range("I132").Value = 6538.7144
$ws.Range("K132").Value = 19616.1432
$ws.Range("M132").Value = -17086.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3675.182
$ws.Range("I22").Value = 5118
$ws.Range("J22").Value = 1150.25
$ws.Range("K22").Value = 5118
$ws.Range("L22").Value = 1150.25
$ws.Range("M22").Value = -4823
$ws.Range("N22").Value = -1740.25
$ws.Range("H27").Value = 3675.182
$ws.Range("I27").Value = 5118
$ws.Range("J27").Value = 1150.25
$ws.Range("K27").Value = 5118
$ws.Range("L27").Value = 1150.25
$ws.Range("M27").Value = -5011
$ws.Range("N27").Value = -1364.25
$ws.Range("H40").Value = 3686.8262
$ws.Range("I40").Value = 3516.7778
$ws.Range("K40").Value = 3516.7778
$ws.Range("M40").Value = -3380.7778
$ws.Range("H82").Value = 1900.4286
$ws.Range("I82").Value = 1883.3334
$ws.Range("J82").Value = 2003
$ws.Range("K82").Value = 1883.3334
$ws.Range("L82").Value = 2003
$ws.Range("M82").Value = -1522.3334
$ws.Range("N82").Value = -2725
$ws.Range("H85").Value = 1900.4286
$ws.Range("I85").Value = 1883.3334
$ws.Range("J85").Value = 2003
$ws.Range("K85").Value = 1883.3334
$ws.Range("L85").Value = 2003
$ws.Range("M85").Value = -635.3334
$ws.Range("N85").Value = -4499
$ws.Range("H93").Value = 1669.75
$ws.Range("I93").Value = 1478.1538
$ws.Range("K93").Value = 1478.1538
$ws.Range("M93").Value = -230.1538
$ws.Range("H122").Value = 1404635.8
$ws.Range("J122").Value = 7688
$ws.Range("L122").Value = 23064
$ws.Range("N122").Value = -27964
$ws.Range("H132").Value = 2255.3333
$ws.Range("I132").Value = 1617.037
$ws.Range("K132").Value = 4851.111
$ws.Range("M132").Value = -2321.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1670.3448
$ws.Range("I113").Value = 1485.3914
$ws.Range("J113").Value = 2379.3333
$ws.Range("K113").Value = 4456.174199999999
$ws.Range("L113").Value = 7137.999899999999
$ws.Range("M113").Value = -2286.174199999999
$ws.Range("N113").Value = -11477.9999
$ws.Range("H132").Value = 903.2222
$ws.Range("I132").Value = 681.6
$ws.Range("J132").Value = 1536.4286
$ws.Range("K132").Value = 2044.8
$ws.Range("L132").Value = 4609.2858
$ws.Range("M132").Value = 485.1999999999998
$ws.Range("N132").Value = -9669.2858
